$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted as row 394; every existing record
# that was previously on row 394 (Berenjena / Vega Central Mapocho de
# Santiago) down through row 411 shifts down by one row, becoming rows
# 395-412. Inserting a whole row handles that shift for every column.
$ws.Rows.Item(394).Insert()

# Populate the newly inserted row 394 with the new record's data.
$ws.Cells.Item(394, 1).Value = 9
$ws.Cells.Item(394, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(394, 3).Value = "Metropolitana"
$ws.Cells.Item(394, 4).Value = 45267
$ws.Cells.Item(394, 5).Value = 13
$ws.Cells.Item(394, 6).Value = 100112001
$ws.Cells.Item(394, 7).Value = "Berenjena"
$ws.Cells.Item(394, 8).Value = "Sin especificar"
$ws.Cells.Item(394, 9).Value = "Primera"
$ws.Cells.Item(394, 10).Value = 160
$ws.Cells.Item(394, 11).Value = 9000
$ws.Cells.Item(394, 12).Value = 10000
$ws.Cells.Item(394, 13).Value = 9500
$ws.Cells.Item(394, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(394, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(394, 16).Value = 190
$ws.Cells.Item(394, 17).Value = 50
$ws.Cells.Item(394, 18).Value = "Hortaliza"
